# Update existing headers/values and add new "EXPORTER" column with
# port of loading / exporter details.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing columns: rename "CARTONS" -> "TO PAY" and update its value,
# and "GROSS WEIGHT" -> "PORT OF LOADING" and update its value.
$ws.Range("B1").Value = "TO PAY"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2006.92"
$ws.Range("B2").ClearFormats()

$ws.Range("C1").Value = "PORT OF LOADING"
$ws.Range("C2").Value = "Chittagong"

# New column D: EXPORTER
$ws.Range("D1").Value = "EXPORTER"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Borders.LineStyle = 1

$ws.Range("D2").Value = "Fakir Fashion Ltd`n89, Motijheel C/A,`nLucky Chamber (2nd Floor)`nDhaka-1000,Bangladesh.`nE-mail:akon@fakirfashion.com`nFactory:Dohargaon,Baliapara,Rupgonj,`nNarayangonj-1400,Bangladesh.`nGB No.GB801930315"
$ws.Rows.Item(2).AutoFit()
